# Refresh the crypto price/volume table (Price = column D, Volume(1h) = column E)
# with the latest scraped values from the GitHub Actions run.
#
# Column D ("Price") cells are plain text in the source data (European-style
# thousands separators like "70.947.45" or trailing-dot decimals like
# "705.97" that must NOT be reinterpreted as numbers). Setting NumberFormat
# to "@" (Text) before assigning the Value keeps Excel from auto-converting
# number-looking strings to floating point, and resetting Style back to
# "Normal" afterwards avoids leaving a stray text-format style behind so the
# cell formatting matches the original file.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.947.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.850.92"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "705.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.849.18"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.51%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -0.69%  "
$ws.Range("E10").Value = "  -0.33%  "
$ws.Range("E11").Value = "  -1.09%  "
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000257"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.499.97"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.830.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.82%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "71.014.58"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("E19").Value = "  +0.71%  "
$ws.Range("E20").Value = "  -2.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "492.97"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.78%  "
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.10%  "
$ws.Range("E25").Value = "  +1.85%  "
$ws.Range("E26").Value = "  +1.28%  "
$ws.Range("E27").Value = "  -2.12%  "
$ws.Range("E28").Value = "  -1.88%  "
$ws.Range("E29").Value = "  +3.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("E32").Value = "  -0.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.29%  "
$ws.Range("E34").Value = "  +0.13%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.47%  "
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.806.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.57%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("E38").Value = "  +0.66%  "
$ws.Range("E39").Value = "  +6.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.07"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.42%  "
$ws.Range("E41").Value = "  +6.57%  "
$ws.Range("E42").Value = "  -4.86%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "163.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("E46").Value = "  -5.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "48.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.69%  "
$ws.Range("E48").Value = "  +0.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "414.73"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.14%  "
$ws.Range("E50").Value = "  -0.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.90%  "
